$d = $word.ActiveDocument

$d.Content.Find.Execute("17×93=", $true, $false, $false, $false, $false, $true, 1, $false, "58×45=", 2) | Out-Null
$d.Content.Find.Execute("47×82=", $true, $false, $false, $false, $false, $true, 1, $false, "69×83=", 2) | Out-Null
$d.Content.Find.Execute("89×91=", $true, $false, $false, $false, $false, $true, 1, $false, "24×56=", 2) | Out-Null
$d.Content.Find.Execute("65×67=", $true, $false, $false, $false, $false, $true, 1, $false, "81×71=", 2) | Out-Null
$d.Content.Find.Execute("23×18=", $true, $false, $false, $false, $false, $true, 1, $false, "98×93=", 2) | Out-Null
$d.Content.Find.Execute("33×19=", $true, $false, $false, $false, $false, $true, 1, $false, "28×88=", 2) | Out-Null
$d.Content.Find.Execute("25×25=", $true, $false, $false, $false, $false, $true, 1, $false, "52×58=", 2) | Out-Null
$d.Content.Find.Execute("75×45=", $true, $false, $false, $false, $false, $true, 1, $false, "61×94=", 2) | Out-Null
$d.Content.Find.Execute("33×97=", $true, $false, $false, $false, $false, $true, 1, $false, "19×29=", 2) | Out-Null
$d.Content.Find.Execute("27×55=", $true, $false, $false, $false, $false, $true, 1, $false, "44×83=", 2) | Out-Null
$d.Content.Find.Execute("35×38=", $true, $false, $false, $false, $false, $true, 1, $false, "34×44=", 2) | Out-Null
$d.Content.Find.Execute("13×67=", $true, $false, $false, $false, $false, $true, 1, $false, "48×86=", 2) | Out-Null
$d.Content.Find.Execute("13×81=", $true, $false, $false, $false, $false, $true, 1, $false, "11×75=", 2) | Out-Null
$d.Content.Find.Execute("51×21=", $true, $false, $false, $false, $false, $true, 1, $false, "15×63=", 2) | Out-Null
$d.Content.Find.Execute("41×65=", $true, $false, $false, $false, $false, $true, 1, $false, "51×43=", 2) | Out-Null
$d.Content.Find.Execute("60×82=", $true, $false, $false, $false, $false, $true, 1, $false, "47×34=", 2) | Out-Null
$d.Content.Find.Execute("39×48=", $true, $false, $false, $false, $false, $true, 1, $false, "91×78=", 2) | Out-Null
$d.Content.Find.Execute("35×74=", $true, $false, $false, $false, $false, $true, 1, $false, "26×48=", 2) | Out-Null
$d.Content.Find.Execute("65×80=", $true, $false, $false, $false, $false, $true, 1, $false, "13×33=", 2) | Out-Null
$d.Content.Find.Execute("49×36=", $true, $false, $false, $false, $false, $true, 1, $false, "67×54=", 2) | Out-Null
$d.Content.Find.Execute("54×94=", $true, $false, $false, $false, $false, $true, 1, $false, "49×83=", 2) | Out-Null
$d.Content.Find.Execute("82×14=", $true, $false, $false, $false, $false, $true, 1, $false, "36×96=", 2) | Out-Null
$d.Content.Find.Execute("89×16=", $true, $false, $false, $false, $false, $true, 1, $false, "14×15=", 2) | Out-Null
$d.Content.Find.Execute("19×55=", $true, $false, $false, $false, $false, $true, 1, $false, "81×39=", 2) | Out-Null
$d.Content.Find.Execute("68×97=", $true, $false, $false, $false, $false, $true, 1, $false, "46×75=", 2) | Out-Null
